$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ccl2"
$ws.Cells.Item(2, 3).Value = "Ackr4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 9.380719000000001
$ws.Cells.Item(2, 8).Value = 28.142157
$ws.Cells.Item(2, 9).Value = 0.03679977590837273
$ws.Cells.Item(2, 10).Value = 0.03679977590837273
$ws.Cells.Item(2, 11).Value = 1.0
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.377371
$ws.Cells.Item(2, 14).Value = 1.132113
$ws.Cells.Item(2, 15).Value = 0.4698794580655765
$ws.Cells.Item(2, 16).Value = 0.4698794580655764
$ws.Cells.Item(2, 17).Value = 3.540011309749
$ws.Cells.Item(2, 18).Value = 31.860101787741
$ws.Cells.Item(2, 19).Value = 0.01729145876076084
$ws.Cells.Item(2, 20).Value = 0.01729145876076083

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ccl2"
$ws.Cells.Item(3, 3).Value = "Ackr4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 9.380719000000001
$ws.Cells.Item(3, 8).Value = 28.142157
$ws.Cells.Item(3, 9).Value = 0.03679977590837273
$ws.Cells.Item(3, 10).Value = 0.03679977590837273
$ws.Cells.Item(3, 11).Value = 1.0
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.3560133333333333
$ws.Cells.Item(3, 14).Value = 1.06804
$ws.Cells.Item(3, 15).Value = 0.443286188209444
$ws.Cells.Item(3, 16).Value = 0.443286188209444
$ws.Cells.Item(3, 17).Value = 3.339661040253334
$ws.Cells.Item(3, 18).Value = 30.05694936228
$ws.Cells.Item(3, 19).Value = 0.01631283238938428
$ws.Cells.Item(3, 20).Value = 0.01631283238938428

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ccl2"
$ws.Cells.Item(4, 3).Value = "Ackr4"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 9.380719000000001
$ws.Cells.Item(4, 8).Value = 28.142157
$ws.Cells.Item(4, 9).Value = 0.03679977590837273
$ws.Cells.Item(4, 10).Value = 0.03679977590837273
$ws.Cells.Item(4, 11).Value = 2.0
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.06973866666666667
$ws.Cells.Item(4, 14).Value = 0.209216
$ws.Cells.Item(4, 15).Value = 0.08683435372497944
$ws.Cells.Item(4, 16).Value = 0.08683435372497944
$ws.Cells.Item(4, 17).Value = 0.6541988354346667
$ws.Cells.Item(4, 18).Value = 5.887789518912
$ws.Cells.Item(4, 19).Value = 0.003195484758227615
$ws.Cells.Item(4, 20).Value = 0.003195484758227615

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ccl2"
$ws.Cells.Item(5, 3).Value = "Ackr4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 90.439374
$ws.Cells.Item(5, 8).Value = 271.318122
$ws.Cells.Item(5, 9).Value = 0.3547860986448385
$ws.Cells.Item(5, 10).Value = 0.3547860986448385
$ws.Cells.Item(5, 11).Value = 1.0
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.377371
$ws.Cells.Item(5, 14).Value = 1.132113
$ws.Cells.Item(5, 15).Value = 0.4698794580655765
$ws.Cells.Item(5, 16).Value = 0.4698794580655764
$ws.Cells.Item(5, 17).Value = 34.12919700575399
$ws.Cells.Item(5, 18).Value = 307.162773051786
$ws.Cells.Item(5, 19).Value = 0.1667066997604369
$ws.Cells.Item(5, 20).Value = 0.1667066997604368

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ccl2"
$ws.Cells.Item(6, 3).Value = "Ackr4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 90.439374
$ws.Cells.Item(6, 8).Value = 271.318122
$ws.Cells.Item(6, 9).Value = 0.3547860986448385
$ws.Cells.Item(6, 10).Value = 0.3547860986448385
$ws.Cells.Item(6, 11).Value = 1.0
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.3560133333333333
$ws.Cells.Item(6, 14).Value = 1.06804
$ws.Cells.Item(6, 15).Value = 0.443286188209444
$ws.Cells.Item(6, 16).Value = 0.443286188209444
$ws.Cells.Item(6, 17).Value = 32.19762300232
$ws.Cells.Item(6, 18).Value = 289.77860702088
$ws.Cells.Item(6, 19).Value = 0.1572717772979703
$ws.Cells.Item(6, 20).Value = 0.1572717772979703

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ccl2"
$ws.Cells.Item(7, 3).Value = "Ackr4"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 90.439374
$ws.Cells.Item(7, 8).Value = 271.318122
$ws.Cells.Item(7, 9).Value = 0.3547860986448385
$ws.Cells.Item(7, 10).Value = 0.3547860986448385
$ws.Cells.Item(7, 11).Value = 2.0
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.06973866666666667
$ws.Cells.Item(7, 14).Value = 0.209216
$ws.Cells.Item(7, 15).Value = 0.08683435372497944
$ws.Cells.Item(7, 16).Value = 0.08683435372497944
$ws.Cells.Item(7, 17).Value = 6.307121356928
$ws.Cells.Item(7, 18).Value = 56.764092212352
$ws.Cells.Item(7, 19).Value = 0.03080762158643136
$ws.Cells.Item(7, 20).Value = 0.03080762158643136

$ws.Cells.Item(8, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value = "Ccl2"
$ws.Cells.Item(8, 3).Value = "Ackr4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 100.179423
$ws.Cells.Item(8, 8).Value = 300.538269
$ws.Cells.Item(8, 9).Value = 0.3929954960840508
$ws.Cells.Item(8, 10).Value = 0.3929954960840508
$ws.Cells.Item(8, 11).Value = 1.0
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.377371
$ws.Cells.Item(8, 14).Value = 1.132113
$ws.Cells.Item(8, 15).Value = 0.4698794580655765
$ws.Cells.Item(8, 16).Value = 0.4698794580655764
$ws.Cells.Item(8, 17).Value = 37.80480903693299
$ws.Cells.Item(8, 18).Value = 340.243281332397
$ws.Cells.Item(8, 19).Value = 0.1846605107221861
$ws.Cells.Item(8, 20).Value = 0.1846605107221861

$ws.Cells.Item(9, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value = "Ccl2"
$ws.Cells.Item(9, 3).Value = "Ackr4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 100.179423
$ws.Cells.Item(9, 8).Value = 300.538269
$ws.Cells.Item(9, 9).Value = 0.3929954960840508
$ws.Cells.Item(9, 10).Value = 0.3929954960840508
$ws.Cells.Item(9, 11).Value = 1.0
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.3560133333333333
$ws.Cells.Item(9, 14).Value = 1.06804
$ws.Cells.Item(9, 15).Value = 0.443286188209444
$ws.Cells.Item(9, 16).Value = 0.443286188209444
$ws.Cells.Item(9, 17).Value = 35.66521031364
$ws.Cells.Item(9, 18).Value = 320.9868928227601
$ws.Cells.Item(9, 19).Value = 0.1742094754425784
$ws.Cells.Item(9, 20).Value = 0.1742094754425784

$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Ccl2"
$ws.Cells.Item(10, 3).Value = "Ackr4"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 100.179423
$ws.Cells.Item(10, 8).Value = 300.538269
$ws.Cells.Item(10, 9).Value = 0.3929954960840508
$ws.Cells.Item(10, 10).Value = 0.3929954960840508
$ws.Cells.Item(10, 11).Value = 2.0
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.06973866666666667
$ws.Cells.Item(10, 14).Value = 0.209216
$ws.Cells.Item(10, 15).Value = 0.08683435372497944
$ws.Cells.Item(10, 16).Value = 0.08683435372497944
$ws.Cells.Item(10, 17).Value = 6.986379387456
$ws.Cells.Item(10, 18).Value = 62.87741448710401
$ws.Cells.Item(10, 19).Value = 0.03412550991928624
$ws.Cells.Item(10, 20).Value = 0.03412550991928624

$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Ccl2"
$ws.Cells.Item(11, 3).Value = "Ackr4"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 1.427630666666667
$ws.Cells.Item(11, 8).Value = 4.282892
$ws.Cells.Item(11, 9).Value = 0.005600475679236752
$ws.Cells.Item(11, 10).Value = 0.005600475679236752
$ws.Cells.Item(11, 11).Value = 1.0
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.377371
$ws.Cells.Item(11, 14).Value = 1.132113
$ws.Cells.Item(11, 15).Value = 0.4698794580655765
$ws.Cells.Item(11, 16).Value = 0.4698794580655764
$ws.Cells.Item(11, 17).Value = 0.5387464123106667
$ws.Cells.Item(11, 18).Value = 4.848717710796
$ws.Cells.Item(11, 19).Value = 0.002631548477069206
$ws.Cells.Item(11, 20).Value = 0.002631548477069206

$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Ccl2"
$ws.Cells.Item(12, 3).Value = "Ackr4"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 1.427630666666667
$ws.Cells.Item(12, 8).Value = 4.282892
$ws.Cells.Item(12, 9).Value = 0.005600475679236752
$ws.Cells.Item(12, 10).Value = 0.005600475679236752
$ws.Cells.Item(12, 11).Value = 1.0
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.3560133333333333
$ws.Cells.Item(12, 14).Value = 1.06804
$ws.Cells.Item(12, 15).Value = 0.443286188209444
$ws.Cells.Item(12, 16).Value = 0.443286188209444
$ws.Cells.Item(12, 17).Value = 0.508255552408889
$ws.Cells.Item(12, 18).Value = 4.574299971680001
$ws.Cells.Item(12, 19).Value = 0.002482613516008556
$ws.Cells.Item(12, 20).Value = 0.002482613516008556

$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Ccl2"
$ws.Cells.Item(13, 3).Value = "Ackr4"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 1.427630666666667
$ws.Cells.Item(13, 8).Value = 4.282892
$ws.Cells.Item(13, 9).Value = 0.005600475679236752
$ws.Cells.Item(13, 10).Value = 0.005600475679236752
$ws.Cells.Item(13, 11).Value = 2.0
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.06973866666666667
$ws.Cells.Item(13, 14).Value = 0.209216
$ws.Cells.Item(13, 15).Value = 0.08683435372497944
$ws.Cells.Item(13, 16).Value = 0.08683435372497944
$ws.Cells.Item(13, 17).Value = 0.09956105918577779
$ws.Cells.Item(13, 18).Value = 0.8960495326720002
$ws.Cells.Item(13, 19).Value = 0.0004863136861589886
$ws.Cells.Item(13, 20).Value = 0.0004863136861589886

$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Ccl2"
$ws.Cells.Item(14, 3).Value = "Ackr4"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3.0
$ws.Cells.Item(14, 6).Value = 1.0
$ws.Cells.Item(14, 7).Value = 53.48524799999999
$ws.Cells.Item(14, 8).Value = 160.455744
$ws.Cells.Item(14, 9).Value = 0.2098181536835013
$ws.Cells.Item(14, 10).Value = 0.2098181536835013
$ws.Cells.Item(14, 11).Value = 1.0
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.377371
$ws.Cells.Item(14, 14).Value = 1.132113
$ws.Cells.Item(14, 15).Value = 0.4698794580655765
$ws.Cells.Item(14, 16).Value = 0.4698794580655764
$ws.Cells.Item(14, 17).Value = 20.18378152300799
$ws.Cells.Item(14, 18).Value = 181.654033707072
$ws.Cells.Item(14, 19).Value = 0.09858924034512341
$ws.Cells.Item(14, 20).Value = 0.0985892403451234

$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Ccl2"
$ws.Cells.Item(15, 3).Value = "Ackr4"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3.0
$ws.Cells.Item(15, 6).Value = 1.0
$ws.Cells.Item(15, 7).Value = 53.48524799999999
$ws.Cells.Item(15, 8).Value = 160.455744
$ws.Cells.Item(15, 9).Value = 0.2098181536835013
$ws.Cells.Item(15, 10).Value = 0.2098181536835013
$ws.Cells.Item(15, 11).Value = 1.0
$ws.Cells.Item(15, 12).Value = 0.3333333333333333
$ws.Cells.Item(15, 13).Value = 0.3560133333333333
$ws.Cells.Item(15, 14).Value = 1.06804
$ws.Cells.Item(15, 15).Value = 0.443286188209444
$ws.Cells.Item(15, 16).Value = 0.443286188209444
$ws.Cells.Item(15, 17).Value = 19.04146142464
$ws.Cells.Item(15, 18).Value = 171.37315282176
$ws.Cells.Item(15, 19).Value = 0.09300948956350259
$ws.Cells.Item(15, 20).Value = 0.09300948956350259

$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Ccl2"
$ws.Cells.Item(16, 3).Value = "Ackr4"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 3.0
$ws.Cells.Item(16, 6).Value = 1.0
$ws.Cells.Item(16, 7).Value = 53.48524799999999
$ws.Cells.Item(16, 8).Value = 160.455744
$ws.Cells.Item(16, 9).Value = 0.2098181536835013
$ws.Cells.Item(16, 10).Value = 0.2098181536835013
$ws.Cells.Item(16, 11).Value = 2.0
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.06973866666666667
$ws.Cells.Item(16, 14).Value = 0.209216
$ws.Cells.Item(16, 15).Value = 0.08683435372497944
$ws.Cells.Item(16, 16).Value = 0.08683435372497944
$ws.Cells.Item(16, 17).Value = 3.729989881856
$ws.Cells.Item(16, 18).Value = 33.569908936704
$ws.Cells.Item(16, 19).Value = 0.01821942377487525
$ws.Cells.Item(16, 20).Value = 0.01821942377487525

$ws.Range("A17:T21").Delete()
